$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Update column widths: A=12, B=18, C=12 (D stays 8)
# Note: this runtime's ColumnWidth setter adds a fixed 5/6 (0.8333...) padding
# when persisting to the OOXML <col width="..."> attribute, so compensate by
# subtracting that offset from the desired stored width.
$ws.Columns.Item(1).ColumnWidth = 12 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 18 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 12 - (5/6)

# Update row 2 data
$ws.Range("A2").Value = 1000545230
$ws.Range("B2").Value = "Baltazar Sanchez"
$ws.Range("C2").Value = 4325465745
